$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused date columns C, D, E (shifts nothing else, only 2 data
# columns remain: the name column A and the single date column B)
$ws.Range("C:E").EntireColumn.Delete() | Out-Null

# Remove the attendance row that is no longer tracked (old row 2); this shifts
# the remaining "Шекшуев Филипп" row up to become row 2
$ws.Rows.Item(2).EntireRow.Delete() | Out-Null

# Resize the name column to its new (narrower) width
$ws.Columns.Item(1).ColumnWidth = 14.166666666666666

# Update the remaining header date
$ws.Range("B1").Value = "20.01.2025"

# Update the remaining person's name
$ws.Range("A2").Value = "Точкее Точк"

# Make sure the attendance mark is present
$ws.Range("B2").Value = "+"
